# Reorders each calendar-year block of monthly rows so that the
# October/November/December rows move to the front of their year,
# followed by January-September (matches the target row layout).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, month label (A), then the four index values (B:E)
$newRowData = @(
  @(2, "2014-10", 97.51690000000001, 100.3939, 103.1357, 81.9134),
  @(3, "2014-11", 96.33280000000001, 100.1969, 103.454, 80.6784),
  @(4, "2014-12", 96.1832, 100.2139, 103.1525, 79.70359999999999),
  @(5, "2014-01", 93.68049999999999, 101.5358, 101.3127, 100.8481),
  @(6, "2014-02", 94.16289999999999, 101.7553, 101.2994, 100.4454),
  @(7, "2014-03", 94.21980000000001, 101.3523, 101.1377, 97.8017),
  @(8, "2014-04", 94.6366, 101.1075, 101.0526, 95.1105),
  @(9, "2014-05", 95.07729999999999, 100.9209, 101.1065, 90.88639999999999),
  @(10, "2014-06", 94.0638, 100.9668, 101.2768, 90.9965),
  @(11, "2014-07", 94.3241, 100.6741, 101.0236, 90.37569999999999),
  @(12, "2014-08", 95.03100000000001, 100.881, 102.2023, 87.5872),
  @(13, "2014-09", 95.9667, 100.3221, 102.7208, 84.42100000000001),
  @(14, "2015-10", 96.3, 99.2, 98.09999999999999, 85.3),
  @(15, "2015-11", 96.6434, 99.2817, 97.6264, 85.988),
  @(16, "2015-12", 96.5565, 98.9748, 97.2046, 86.2958),
  @(17, "2015-01", 96.3357, 100.2694, 102.9629, 79.18000000000001),
  @(18, "2015-02", 95.8858, 100.0777, 102.5267, 78.253),
  @(19, "2015-03", 95.5891, 100.2298, 101.7725, 76.5829),
  @(20, "2015-04", 96.01130000000001, 99.6118, 102.0081, 76.8193),
  @(21, "2015-05", 95.9269, 99.6125, 101.0954, 77.4044),
  @(22, "2015-06", 96.88290000000001, 99.69370000000001, 100.7839, 78.14579999999999),
  @(23, "2015-07", 98.19289999999999, 99.97320000000001, 100.4273, 78.6752),
  @(24, "2015-08", 97.02200000000001, 99.4611, 98.9436, 80.3237),
  @(25, "2015-09", 96.96080000000001, 99.7175, 98.57989999999999, 83.1073),
  @(26, "2016-10", 96, 100.7, 98.2, 99.7),
  @(27, "2016-11", 96.59999999999999, 101, 98.5, 104.9),
  @(28, "2016-12", 96.90000000000001, 101.4, 99.09999999999999, 107.7),
  @(29, "2016-01", 97.6063, 99.5716, 98.3573, 89.21040000000001),
  @(30, "2016-02", 96.9828, 99.00539999999999, 98.4607, 89.73990000000001),
  @(31, "2016-03", 96.4183, 99.0347, 98.065, 90.55800000000001),
  @(32, "2016-04", 96.15689999999999, 99.5594, 97.6434, 91.35890000000001),
  @(33, "2016-05", 96.7, 99.90000000000001, 97.90000000000001, 92.2),
  @(34, "2016-06", 97.40000000000001, 99.90000000000001, 98, 92.8),
  @(35, "2016-07", 97, 100, 97.5, 94),
  @(36, "2016-08", 97.3, 100, 97.40000000000001, 94.2),
  @(37, "2016-09", 96.3, 100.2, 98, 95.59999999999999),
  @(38, "2017-10", 100.4, 107, 106.5, 119.2),
  @(39, "2017-11", 101.4, 107.5, 107, 115.2),
  @(40, "2017-12", 102.3, 108.3, 107.7, 114),
  @(41, "2017-01", 97.8, 101.4, 100.4, 111),
  @(42, "2017-02", 99.3, 102.4, 100.7, 116.7),
  @(43, "2017-03", 99.3, 102.5, 101.5, 118.6),
  @(44, "2017-04", 99.8, 102.8, 102.1, 119),
  @(45, "2017-05", 99, 103.1, 102.7, 120.1),
  @(46, "2017-06", 98.7, 104.4, 103.3, 119.2),
  @(47, "2017-07", 98.3, 104.9, 105, 118.7),
  @(48, "2017-08", 98.59999999999999, 105.4, 106.4, 118.3),
  @(49, "2017-09", 99.90000000000001, 106.1, 106.3, 117.8)
)

foreach ($entry in $newRowData) {
  $r = $entry[0]
  $ws.Cells.Item($r, 1).Value = $entry[1]
  $ws.Cells.Item($r, 2).Value = $entry[2]
  $ws.Cells.Item($r, 3).Value = $entry[3]
  $ws.Cells.Item($r, 4).Value = $entry[4]
  $ws.Cells.Item($r, 5).Value = $entry[5]
}
